# Apply updated values to column F ("dSF") on Sheet1.
# Mapping of row -> new value, per the commit's "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -5
    4  = -7
    7  = -3
    9  = -2
    10 = -2
    13 = 5
    17 = 2
    18 = -1
    19 = -1
    23 = 0
    24 = 2
    25 = 6
    27 = -5
    28 = -3
    29 = -3
    30 = -2
    32 = -1
    35 = 0
    38 = 0
    43 = -5
    55 = -6
    56 = -6
    57 = 5
    58 = 2
    61 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
